# Update the cryptocurrency price table (rows 2-51) with the latest scraped
# values. A new coin (LEO) enters the ranking at row 26, pushing the coins
# that were ranked 26-50 down by one row each; the former last row (Aave,
# row 51) drops off the bottom of the table and is replaced by Quant.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.999.07'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '1.863.84'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = "'311.91"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = "'0.5124"
$ws.Range("E7").Value = '  +1.95%  '
$ws.Range("D8").Value = "'0.3831"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = "'0.08231"
$ws.Range("E9").Value = '  -4.10%  '
$ws.Range("D10").Value = "'1.110"
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").Value = "'41.39"
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = "'6.187"
$ws.Range("E12").Value = '  -1.73%  '
$ws.Range("D13").Value = "'20.49"
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").Value = '1.853.70'
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = "'7.279"
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = "'0.00001094"
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = "'90.73"
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").Value = "'0.06631"
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = "'17.63"
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = "'6.016"
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").Value = '28.035.65'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").Value = "'11.05"
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").Value = "'2.246"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = "'3.393"
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '2.075.89'
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'2.463"
$ws.Range("E28").Value = '  -5.40%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = "'157.68"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = "'20.50"
$ws.Range("E30").Value = '  -1.08%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = "'124.74"
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = "'0.1063"
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = "'1.027"
$ws.Range("E33").Value = '  -2.81%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = "'5.828"
$ws.Range("E34").Value = '  +3.70%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = "'3.601"
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = "'9.448"
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.02420"
$ws.Range("E37").Value = '  -1.36%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = "'0.06475"
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = "'0.2192"
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.6488"
$ws.Range("E40").Value = '  +1.45%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = "'1.191"
$ws.Range("E41").Value = '  -2.22%  '
$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").Value = "'4.984"
$ws.Range("E42").Value = '  +2.00%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = "'1.210"
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = "'11.21"
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.6097"
$ws.Range("E45").Value = '  +1.29%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'13.11"
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = "'1.285"
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = "'3.667"
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'2.014"
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = "'1.219"
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = "'120.83"
$ws.Range("E51").Value = '  +0.17%  '
